$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "69.503.73"
$ws.Cells.Item(2, 5).Value = "  -0.44%  "
$ws.Cells.Item(3, 4).Value = "3.776.07"
$ws.Cells.Item(3, 5).Value = "  -0.03%  "
$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.12%  "
$ws.Cells.Item(5, 4).Value = "'614.04"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.91%  "
$ws.Cells.Item(6, 4).Value = "'177.78"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +0.36%  "
$ws.Cells.Item(7, 4).Value = "3.773.69"
$ws.Cells.Item(7, 5).Value = "  -0.10%  "
$ws.Cells.Item(8, 5).Value = "  +0.12%  "
$ws.Cells.Item(9, 4).Value = "'0.529"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -1.06%  "
$ws.Cells.Item(10, 5).Value = "  -1.78%  "
$ws.Cells.Item(11, 5).Value = "  +4.52%  "
$ws.Cells.Item(12, 5).Value = "  -0.89%  "
$ws.Cells.Item(13, 4).Value = "'39.89"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -2.06%  "
$ws.Cells.Item(14, 4).Value = "'0.0000254"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -1.96%  "
$ws.Cells.Item(15, 4).Value = "4.403.73"
$ws.Cells.Item(15, 5).Value = "  +0.13%  "
$ws.Cells.Item(16, 4).Value = "3.769.17"
$ws.Cells.Item(16, 5).Value = "  -0.28%  "
$ws.Cells.Item(17, 4).Value = "69.535.75"
$ws.Cells.Item(17, 5).Value = "  -0.47%  "
$ws.Cells.Item(18, 4).Value = "'7.58"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +0.25%  "
$ws.Cells.Item(19, 5).Value = "  -3.26%  "
$ws.Cells.Item(20, 4).Value = "'510.34"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +0.49%  "
$ws.Cells.Item(21, 4).Value = "'16.35"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -2.47%  "
$ws.Cells.Item(22, 5).Value = "  -2.19%  "
$ws.Cells.Item(23, 4).Value = "'0.733"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +1.10%  "
$ws.Cells.Item(24, 5).Value = "  -0.38%  "
$ws.Cells.Item(25, 4).Value = "'86.53"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.64%  "
$ws.Cells.Item(26, 4).Value = "'12.88"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -2.16%  "
$ws.Cells.Item(27, 4).Value = "'0.0000136"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -2.06%  "
$ws.Cells.Item(28, 5).Value = "  -4.02%  "
$ws.Cells.Item(29, 5).Value = "  +0.18%  "
$ws.Cells.Item(30, 5).Value = "  +2.06%  "
$ws.Cells.Item(31, 5).Value = "  +2.68%  "
$ws.Cells.Item(32, 4).Value = "'8.10"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +3.90%  "
$ws.Cells.Item(33, 4).Value = "'30.76"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -1.07%  "
$ws.Cells.Item(34, 5).Value = "  +0.01%  "
$ws.Cells.Item(35, 4).Value = "'0.998"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -0.18%  "
$ws.Cells.Item(36, 5).Value = "  -2.38%  "
$ws.Cells.Item(37, 4).Value = "'6.13"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -0.78%  "
$ws.Cells.Item(38, 5).Value = "  +6.43%  "
$ws.Cells.Item(39, 4).Value = "'0.341"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +2.30%  "
$ws.Cells.Item(40, 4).Value = "'459.04"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +8.60%  "
$ws.Cells.Item(41, 4).Value = "'2.08"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -1.51%  "
$ws.Cells.Item(42, 4).Value = "'49.82"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -2.07%  "
$ws.Cells.Item(43, 4).Value = "'2.98"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +6.79%  "
$ws.Cells.Item(44, 4).Value = "'44.45"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -2.26%  "
$ws.Cells.Item(45, 5).Value = "  -1.26%  "
$ws.Cells.Item(46, 4).Value = "2.957.84"
$ws.Cells.Item(46, 5).Value = "  -2.52%  "
$ws.Cells.Item(47, 5).Value = "  +0.02%  "
$ws.Cells.Item(48, 5).Value = "  +0.05%  "
$ws.Cells.Item(49, 5).Value = "  -0.01%  "
$ws.Cells.Item(50, 4).Value = "'139.06"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +0.32%  "
$ws.Cells.Item(51, 4).Value = "'2.47"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -0.26%  "
